$d = $word.ActiveDocument

$pairs = @(
    @("939×7=6573", "394×9=3546"),
    @("679×3=2037", "791×9=7119"),
    @("484×9=4356", "195×9=1755"),
    @("539×8=4312", "693×6=4158"),
    @("437×7=3059", "432×4=1728"),
    @("690×2=1380", "565×7=3955"),
    @("593×6=3558", "353×7=2471"),
    @("556×8=4448", "862×5=4310"),
    @("651×7=4557", "813×6=4878"),
    @("215×6=1290", "217×9=1953"),
    @("170×3=510",  "722×7=5054"),
    @("199×6=1194", "962×3=2886"),
    @("294×3=882",  "320×9=2880"),
    @("506×9=4554", "490×5=2450"),
    @("490×4=1960", "235×6=1410"),
    @("231×3=693",  "207×2=414"),
    @("359×2=718",  "479×8=3832"),
    @("770×3=2310", "154×9=1386"),
    @("531×9=4779", "941×9=8469"),
    @("347×4=1388", "520×6=3120"),
    @("264×3=792",  "449×7=3143"),
    @("311×5=1555", "919×6=5514"),
    @("558×7=3906", "947×4=3788"),
    @("589×7=4123", "450×6=2700"),
    @("951×3=2853", "994×6=5964")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
